$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.953.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.07%  "

$ws.Range("D3").Value = "'3.095.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'522.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.88%  "

$ws.Range("D6").Value = "'143.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.05%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("D9").Value = "'7.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").Value = "'0.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("D11").Value = "'0.382"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.82%  "

$ws.Range("D12").Value = "'3.630.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").Value = "'26.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.37%  "

$ws.Range("D15").Value = "'0.0000166"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").Value = "'58.971.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "

$ws.Range("D17").Value = "'3.098.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "'6.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").Value = "'12.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").Value = "'8.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "'343.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'0.506"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.61%  "

$ws.Range("D24").Value = "'65.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").Value = "'0.0₃0924"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").Value = "'6.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.75%  "

$ws.Range("D29").Value = "'7.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.42%  "

$ws.Range("D30").Value = "'1.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("E31").Value = "  +2.94%  "

$ws.Range("D32").Value = "'20.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.35%  "

$ws.Range("D33").Value = "'155.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("D34").Value = "'4.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.61%  "

$ws.Range("D35").Value = "'6.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.36%  "

$ws.Range("D36").Value = "'26.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.51%  "

$ws.Range("E37").Value = "  +4.72%  "

$ws.Range("D38").Value = "'0.0686"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.58%  "

$ws.Range("D39").Value = "'3.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.46%  "

$ws.Range("D40").Value = "'3.140.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("D41").Value = "'36.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.34%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D44").Value = "'1.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.50%  "

$ws.Range("D45").Value = "'2.286.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("E46").Value = "  +1.45%  "

$ws.Range("D47").Value = "'20.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.71%  "

$ws.Range("D48").Value = "'0.962"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").Value = "'6.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("D50").Value = "'0.753"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.35%  "

$ws.Range("D51").Value = "'264.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.99%  "
